$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume number + week-of dates) ---
$ws.Range("A8").Value = "Volume 32   Number  24"
$ws.Range("C9").Value = "Report Covering the Week  6/9/2025  Through  6/15/2025"

# --- Crime-stat grid (rows 14-33): new weekly figures ---
# Row 14
$ws.Range("A14").Value = 'Murder'
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = '0'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '***.*'
$ws.Range("F14").NumberFormat = "@"
$ws.Range("F14").Value = '0'
$ws.Range("G14").Value = 5
$ws.Range("H14").Value = -100
$ws.Range("I14").Value = 13
$ws.Range("J14").Value = 25
$ws.Range("K14").Value = -48
$ws.Range("L14").Value = -48
$ws.Range("M14").Value = -35
$ws.Range("N14").Value = -91.275167785234

# Row 15
$ws.Range("A15").Value = 'Rape'
$ws.Range("C15").Value = 3
$ws.Range("D15").Value = 7
$ws.Range("E15").Value = -57.142857142857
$ws.Range("F15").Value = 12
$ws.Range("G15").Value = 16
$ws.Range("H15").Value = -25
$ws.Range("I15").Value = 87
$ws.Range("J15").Value = 66
$ws.Range("K15").Value = 31.818181818181
$ws.Range("L15").Value = 35.9375
$ws.Range("M15").Value = -5.434782608695
$ws.Range("N15").Value = -61.160714285714

# Row 16
$ws.Range("A16").Value = 'Robbery'
$ws.Range("C16").Value = 31
$ws.Range("D16").Value = 27
$ws.Range("E16").Value = 14.814814814814
$ws.Range("F16").Value = 147
$ws.Range("G16").Value = 175
$ws.Range("H16").Value = -16
$ws.Range("I16").Value = 765
$ws.Range("J16").Value = 946
$ws.Range("K16").Value = -19.133192389006
$ws.Range("L16").Value = -5.788177339901
$ws.Range("M16").Value = -26.158301158301
$ws.Range("N16").Value = -81.597305749338

# Row 17
$ws.Range("A17").Value = 'Fel. Assault'
$ws.Range("C17").Value = 53
$ws.Range("D17").Value = 78
$ws.Range("E17").Value = -32.051282051282
$ws.Range("F17").Value = 231
$ws.Range("G17").Value = 282
$ws.Range("H17").Value = -18.085106382978
$ws.Range("I17").Value = 1295
$ws.Range("J17").Value = 1459
$ws.Range("K17").Value = -11.240575736806
$ws.Range("L17").Value = -2.631578947368
$ws.Range("M17").Value = 56.779661016949
$ws.Range("N17").Value = -48.611111111111

# Row 18
$ws.Range("A18").Value = 'Burglary'
$ws.Range("C18").Value = 19
$ws.Range("D18").Value = 28
$ws.Range("E18").Value = -32.142857142857
$ws.Range("F18").Value = 94
$ws.Range("G18").Value = 103
$ws.Range("H18").Value = -8.737864077669
$ws.Range("I18").Value = 638
$ws.Range("J18").Value = 634
$ws.Range("K18").Value = 0.630914826498
$ws.Range("L18").Value = -11.511789181692
$ws.Range("M18").Value = 7.588532883642
$ws.Range("N18").Value = -86.86431953881

# Row 19
$ws.Range("A19").Value = 'Gr. Larceny'
$ws.Range("C19").Value = 124
$ws.Range("D19").Value = 120
$ws.Range("E19").Value = 3.333333333333
$ws.Range("F19").Value = 499
$ws.Range("G19").Value = 515
$ws.Range("H19").Value = -3.106796116504
$ws.Range("I19").Value = 2679
$ws.Range("J19").Value = 2881
$ws.Range("K19").Value = -7.011454356126
$ws.Range("L19").Value = -2.154857560262
$ws.Range("M19").Value = 30.364963503649
$ws.Range("N19").Value = -44.955824943497

# Row 20
$ws.Range("A20").Value = 'G.L.A.'
$ws.Range("C20").Value = 14
$ws.Range("D20").Value = 19
$ws.Range("E20").Value = -26.315789473684
$ws.Range("F20").Value = 71
$ws.Range("G20").Value = 74
$ws.Range("H20").Value = -4.054054054054
$ws.Range("I20").Value = 383
$ws.Range("J20").Value = 406
$ws.Range("K20").Value = -5.665024630541
$ws.Range("L20").Value = -34.417808219178
$ws.Range("M20").Value = 58.921161825726
$ws.Range("N20").Value = -91.166974169741

# Row 21
$ws.Range("A21").Value = 'TOTAL'
$ws.Range("C21").Value = 244
$ws.Range("D21").Value = 279
$ws.Range("E21").Value = -12.544802867383
$ws.Range("F21").Value = 1054
$ws.Range("G21").Value = 1170
$ws.Range("H21").Value = -9.914529914529
$ws.Range("I21").Value = 5860
$ws.Range("J21").Value = 6417
$ws.Range("K21").Value = -8.680068567866
$ws.Range("L21").Value = -6.598661141217
$ws.Range("M21").Value = 20.501747892247
$ws.Range("N21").Value = -72.240644244433

# Row 22
$ws.Range("A22").Value = 'Transit'
$ws.Range("C22").Value = 2
$ws.Range("D22").Value = 7
$ws.Range("E22").Value = -71.428571428571
$ws.Range("F22").Value = 11
$ws.Range("G22").Value = 15
$ws.Range("H22").Value = -26.666666666666
$ws.Range("I22").Value = 96
$ws.Range("J22").Value = 117
$ws.Range("K22").Value = -17.948717948717
$ws.Range("L22").Value = -28.888888888888
$ws.Range("M22").Value = -9.43396226415
$ws.Range("N22").Value = '***.*'

# Row 23
$ws.Range("C23").Value = 22
$ws.Range("D23").Value = 26
$ws.Range("E23").Value = -15.384615384615
$ws.Range("F23").Value = 84
$ws.Range("G23").Value = 117
$ws.Range("H23").Value = -28.205128205128
$ws.Range("I23").Value = 567
$ws.Range("J23").Value = 607
$ws.Range("K23").Value = -6.58978583196
$ws.Range("L23").Value = 2.53164556962
$ws.Range("M23").Value = 60.623229461756
$ws.Range("N23").Value = '***.*'

# Row 24
$ws.Range("C24").Value = 260
$ws.Range("D24").Value = 242
$ws.Range("E24").Value = 7.438016528925
$ws.Range("F24").Value = 1021
$ws.Range("G24").Value = 1009
$ws.Range("H24").Value = 1.189296333002
$ws.Range("I24").Value = 6301
$ws.Range("J24").Value = 5672
$ws.Range("K24").Value = 11.089562764457
$ws.Range("L24").Value = 0.398342893562
$ws.Range("M24").Value = 60.616874840683
$ws.Range("N24").Value = '***.*'

# Row 25
$ws.Range("C25").Value = 106
$ws.Range("D25").Value = 133
$ws.Range("E25").Value = -20.300751879699
$ws.Range("F25").Value = 484
$ws.Range("G25").Value = 532
$ws.Range("H25").Value = -9.022556390977
$ws.Range("I25").Value = 3406
$ws.Range("J25").Value = 3018
$ws.Range("K25").Value = 12.856196156395
$ws.Range("L25").Value = -4.083356800901
$ws.Range("M25").Value = '***.*'
$ws.Range("N25").Value = '***.*'

# Row 26
$ws.Range("C26").Value = 119
$ws.Range("D26").Value = 113
$ws.Range("E26").Value = 5.309734513274
$ws.Range("F26").Value = 451
$ws.Range("G26").Value = 476
$ws.Range("H26").Value = -5.252100840336
$ws.Range("I26").Value = 2219
$ws.Range("J26").Value = 2283
$ws.Range("K26").Value = -2.803328953131
$ws.Range("L26").Value = 7.561803199224
$ws.Range("M26").Value = -8.267879288962
$ws.Range("N26").Value = '***.*'

# Row 27
$ws.Range("C27").Value = 3
$ws.Range("D27").Value = 8
$ws.Range("E27").Value = -62.5
$ws.Range("F27").Value = 13
$ws.Range("G27").Value = 24
$ws.Range("H27").Value = -45.833333333333
$ws.Range("I27").Value = 108
$ws.Range("J27").Value = 103
$ws.Range("K27").Value = 4.854368932038
$ws.Range("L27").Value = -12.903225806451
$ws.Range("M27").Value = '***.*'
$ws.Range("N27").Value = '***.*'

# Row 28
$ws.Range("C28").Value = 9
$ws.Range("D28").Value = 12
$ws.Range("E28").Value = -25
$ws.Range("F28").Value = 60
$ws.Range("G28").Value = 43
$ws.Range("H28").Value = 39.53488372093
$ws.Range("I28").Value = 287
$ws.Range("J28").Value = 258
$ws.Range("K28").Value = 11.240310077519
$ws.Range("L28").Value = 8.301886792452
$ws.Range("M28").Value = '***.*'
$ws.Range("N28").Value = '***.*'

# Row 29
$ws.Range("C29").Value = 1
$ws.Range("D29").Value = 5
$ws.Range("E29").Value = -80
$ws.Range("F29").Value = 3
$ws.Range("G29").Value = 16
$ws.Range("H29").Value = -81.25
$ws.Range("I29").Value = 32
$ws.Range("J29").Value = 57
$ws.Range("K29").Value = -43.859649122807
$ws.Range("L29").Value = -51.515151515151
$ws.Range("M29").Value = -63.218390804597
$ws.Range("N29").Value = -90.643274853801

# Row 30
$ws.Range("C30").Value = 1
$ws.Range("D30").Value = 3
$ws.Range("E30").Value = -66.666666666666
$ws.Range("F30").Value = 3
$ws.Range("G30").Value = 14
$ws.Range("H30").Value = -78.571428571428
$ws.Range("I30").Value = 29
$ws.Range("J30").Value = 47
$ws.Range("K30").Value = -38.297872340425
$ws.Range("L30").Value = -52.459016393442
$ws.Range("M30").Value = -62.337662337662
$ws.Range("N30").Value = -90.822784810126

# Row 31
$ws.Range("F31").Value = 4
$ws.Range("G31").Value = 11
$ws.Range("H31").Value = -63.636363636363
$ws.Range("I31").Value = 32
$ws.Range("J31").Value = 55
$ws.Range("K31").Value = -41.818181818181
$ws.Range("L31").Value = -8.571428571428
$ws.Range("M31").Value = '***.*'
$ws.Range("N31").Value = '***.*'

# Row 33
$ws.Range("C33").NumberFormat = "@"
$ws.Range("C33").Value = '0'
$ws.Range("D33").Value = 2
$ws.Range("E33").Value = -100
$ws.Range("F33").Value = 1
$ws.Range("G33").Value = 4
$ws.Range("H33").Value = -75
$ws.Range("I33").Value = 8
$ws.Range("J33").Value = 9
$ws.Range("K33").Value = -11.111111111111
$ws.Range("L33").Value = -27.272727272727
$ws.Range("M33").Value = '***.*'
$ws.Range("N33").Value = '***.*'

# Row 40
$ws.Range("A40").Value = 'Robbery'

# Row 41
$ws.Range("A41").Value = 'Fel. Assault'

# Row 42
$ws.Range("A42").Value = 'Burglary'

# Row 43
$ws.Range("A43").Value = 'Gr. Larceny'

# Row 44
$ws.Range("A44").Value = 'G.L.A.'

# Row 45
$ws.Range("A45").Value = 'TOTAL'

# Row 46
$ws.Range("A46").Value = 'TOTAL'

# --- Historical Perspective block (rows 40-46): label shift only ---
$ws.Range("A40").Value = 'Robbery'
$ws.Range("A41").Value = 'Fel. Assault'
$ws.Range("A42").Value = 'Burglary'
$ws.Range("A43").Value = 'Gr. Larceny'
$ws.Range("A44").Value = 'G.L.A.'
$ws.Range("A45").Value = 'TOTAL'
$ws.Range("A46").Value = 'TOTAL'

# --- Column H width cosmetic bestFit tweak ---
$ws.Range("H1").ColumnWidth = 6.71
